$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.873.29'
$ws.Range("E2").Value = '  +1.62%  '

$ws.Range("D3").Value = '3.154.50'
$ws.Range("E3").Value = '  +3.05%  '

$ws.Range("E4").Value = '  +0.10%  '

$ws.Range("D5").Value = '''571.01'
$ws.Range("E5").Value = '  +2.75%  '

$ws.Range("D6").Value = '''150.52'
$ws.Range("E6").Value = '  +5.30%  '

$ws.Range("D7").Value = '''1.00'
$ws.Range("E7").Value = '  +0.12%  '

$ws.Range("D8").Value = '3.152.73'
$ws.Range("E8").Value = '  +2.96%  '

$ws.Range("D9").Value = '''0.528'
$ws.Range("E9").Value = '  +4.84%  '

$ws.Range("D10").Value = '''0.162'
$ws.Range("E10").Value = '  +4.40%  '

$ws.Range("D11").Value = '''6.19'
$ws.Range("E11").Value = '  +1.20%  '

$ws.Range("D12").Value = '''0.504'
$ws.Range("E12").Value = '  +7.28%  '

$ws.Range("D13").Value = '''0.0000261'
$ws.Range("E13").Value = '  +14.24%  '

$ws.Range("D14").Value = '''38.18'
$ws.Range("E14").Value = '  +9.73%  '

$ws.Range("D15").Value = '3.672.75'
$ws.Range("E15").Value = '  +3.34%  '

$ws.Range("D16").Value = '64.938.01'
$ws.Range("E16").Value = '  +1.75%  '

$ws.Range("D17").Value = '''7.20'
$ws.Range("E17").Value = '  +7.53%  '

$ws.Range("D18").Value = '3.154.30'
$ws.Range("E18").Value = '  +3.43%  '

$ws.Range("E19").Value = '  +0.97%  '

$ws.Range("D20").Value = '''515.64'
$ws.Range("E20").Value = '  +7.72%  '

$ws.Range("D21").Value = '''14.91'
$ws.Range("E21").Value = '  +6.73%  '

$ws.Range("D22").Value = '''0.736'
$ws.Range("E22").Value = '  +9.12%  '

$ws.Range("D23").Value = '''15.46'
$ws.Range("E23").Value = '  +10.03%  '

$ws.Range("E24").Value = '  +4.83%  '

$ws.Range("D25").Value = '''84.99'
$ws.Range("E25").Value = '  +4.94%  '

$ws.Range("D26").Value = '''0.999'
$ws.Range("E26").Value = '  -0.13%  '

$ws.Range("E27").Value = '  +5.24%  '

$ws.Range("D28").Value = '''8.90'
$ws.Range("E28").Value = '  +11.96%  '

$ws.Range("E29").Value = '  +7.33%  '

$ws.Range("D30").Value = '''27.85'
$ws.Range("E30").Value = '  +6.56%  '

$ws.Range("E31").Value = '  +0.13%  '

$ws.Range("D32").Value = '''2.71'
$ws.Range("E32").Value = '  +10.24%  '

$ws.Range("D33").Value = '''1.20'
$ws.Range("E33").Value = '  +4.38%  '

$ws.Range("E34").Value = '  +11.10%  '

$ws.Range("E35").Value = '  +7.86%  '

$ws.Range("D36").Value = '''55.96'
$ws.Range("E36").Value = '  +1.51%  '

$ws.Range("D37").Value = '''483.92'
$ws.Range("E37").Value = '  +10.04%  '

$ws.Range("D38").Value = '''0.0868'
$ws.Range("E38").Value = '  +7.71%  '

$ws.Range("D39").Value = '''0.0424'
$ws.Range("E39").Value = '  +4.47%  '

$ws.Range("D40").Value = '''3.03'
$ws.Range("E40").Value = '  +3.20%  '

$ws.Range("D41").Value = '3.119.78'
$ws.Range("E41").Value = '  +5.68%  '

$ws.Range("D42").Value = '''8.67'
$ws.Range("E42").Value = '  +6.05%  '

$ws.Range("E43").Value = '  +4.99%  '

$ws.Range("D44").Value = '''0.291'
$ws.Range("E44").Value = '  +12.89%  '

$ws.Range("D45").Value = '''2.48'
$ws.Range("E45").Value = '  +16.92%  '

$ws.Range("D46").Value = '''29.59'
$ws.Range("E46").Value = '  +4.65%  '

$ws.Range("D47").Value = '0.0₃0578'
$ws.Range("E47").Value = '  +11.97%  '

$ws.Range("E48").Value = '  -0.09%  '

$ws.Range("E49").Value = '  +3.59%  '

$ws.Range("E50").Value = '  +12.06%  '

$ws.Range("D51").Value = '''121.06'
$ws.Range("E51").Value = '  +3.82%  '
